# src: correcciones de carga masiva y testeo funcionando
#
# Applies the data corrections to TestCargaMasivaColaboraciones:
#  - de-accent the header labels "Fecha de colaboración" / "Forma de colaboración"
#  - de-accent "Rocío" / "Ochoa Charlín" for the DNI 44651389 rows
#  - fix the mismatched name/mail pairing on row 11 (Teal -> Teal 2, Music -> Music2)
#  - drop the stray empty styled cell at J11 (no longer needed)
#  - turn on AutoFilter for the Documento column
#  - set the active selection/page setup the way it was left after testing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header corrections (remove accents) ---
$ws.Range("F1").Value = "Fecha de colaboracion"
$ws.Range("G1").Value = "Forma de colaboracion"

# --- rows 13-20 belong to Rocio Ochoa Charlin; remove the accents there too ---
for ($r = 13; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "Rocio"
    $ws.Cells.Item($r, 4).Value = "Ochoa Charlin"
}

# --- row 11 data fix (distinguish this colaborador from the "Teal"/"Music" one) ---
$ws.Range("C11").Value = "Teal 2"
$ws.Range("D11").Value = "Music2"

# --- drop the stray formatted-but-empty cell that was hanging out past the table ---
$ws.Range("J11").Clear()

# --- enable the filter on the "Documento" column, like after testing the import ---
$names = $ws.Names.Add("_xlnm._FilterDatabase", "=Hoja1!`$B`$1:`$B`$20")
$names.Visible = $false

# --- leave the selection where it ended up after scrolling through the test rows ---
$ws.Activate()
$ws.Range("C18").Select()

# --- page setup touched via Page Setup / Print Preview while testing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit applied"
